$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Vin(max) modified as 250V -> update Vin (B2) and Vin(max) (B10) accordingly,
# along with the other dependent input values changed in this revision.
$ws.Range("B2").Value = 250    # Vin
$ws.Range("B5").Value = 0.7    # Vd
$ws.Range("B8").Value = 0.4    # ∆IL
$ws.Range("B10").Value = 250   # Vin(max)
$ws.Range("B11").Value = 0.5   # ESR

# Restore the cursor/selection position left by the author
$ws.Range("C17").Select()
